$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 held a numeric placeholder cell ("0") with a custom bold+bordered
# style; row 2 held the real question-bank text (as a shared string) with
# the default style. Deleting row 1 removes that placeholder/style usage
# entirely and shifts row 2 (and its default styling) up into row 1.
$ws.Rows(1).Delete()

$newText = @'
questions = [
    {
        "title": "You wrote a snippet of code in JavaScript for the browser. Now you want to port that code to Node.js. However, the browser code uses the window object to store and access global variables, which is not available in Node.js. How can you accommodate this?",
        "ques_type": 2,
        "options": [
            "Use Node.js\u2019s global object instead.",
            "Install the Gecko engine to run Node.js on, instead of V8.",
            "Use the document variable instead of the window variable.",
            "Install webpack to convert the browser code to Node.js code."
        ],
        "score": "Use Node.js\u2019s global object instead."
    },
    {
        "title": "You are working as a developer on a project. You want to install some dependencies that will be useful to you while you are writing the code. However, you do not need these dependencies to be available when the application actually runs in production. How should you specify that the dependency be available for developers but not in production?",
        "ques_type": 2,
        "options": [
            "Delete the developer dependencies from the node_modules/ directory in production.",
            "Switch from npm to yarn.",
            "Use the --save-dev flag when installing the package with npm install.",
            "Install the dependency using pip install."
        ],
        "score": "Use the --save-dev flag when installing the package with npm install."
    },
    {
        "title": "Your team wants to select a framework to use on top of Node.js to facilitate API development. Your colleagues explain that they need a mainstream, minimal, unopinionated framework. Since you are the most experienced developer on the team, they ask you which framework they should use. What framework would best meet these requirements?",
        "ques_type": 2,
        "options": [
            "Deno",
            "React",
            "Nuxt.js",
            "Express"
        ],
        "score": "Express"
    },
    {
        "title": "You are writing a web application with Node.js and want to store users\u2019 passwords. After a meeting, management asks you to protect the passwords from common attacks, such as rainbow tables. True or false: You can protect the passwords from rainbow tables by hashing them with salt.",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "True"
    }
]
'@

$ws.Range("A1").Value = $newText

# Re-fit the row height: assigning a multi-line string auto-expands the
# row; AutoFit restores it to the sheet's default (no explicit ht/
# customHeight survives in the saved XML), matching the original layout.
$ws.Rows(1).AutoFit()
